# Swap columns C (codeforiati:group-name) and D (codeforiati:group-code)
# so that the group-code column comes before the group-name column,
# matching the updated codeforIATI codelists export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRows = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $usedRows; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
